$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values: G2:G4 from 60 -> 600
$ws.Range("G2").Value = 600
$ws.Range("G3").Value = 600
$ws.Range("G4").Value = 600

# Update values: F5:F7 from 0.99 -> 0.2
$ws.Range("F5").Value = 0.2
$ws.Range("F6").Value = 0.2
$ws.Range("F7").Value = 0.2

# Update the active selection to F5:F7 with active cell F5
$ws.Range("F5:F7").Select()
